$wb = $excel.ActiveWorkbook

# --- Sheet "trait" (sheet1.xml): fix up row 2 first ---
$ws1 = $wb.Worksheets.Item("trait")
$ws1.Range("C2").Value = 'radius length'
$ws1.Range("A2").Value = 'in oba'

# --- Sheet "axis" (sheet2.xml): populate the whole mapping table ---
$ws2 = $wb.Worksheets.Item("axis")
$ws2.Range("A1").Value = 'Trait'
$ws2.Range("B1").Value = 'Axis'
$ws2.Range("C1").Value = 'Structure'
$ws2.Range("D1").Value = 'Pattern name'
$ws2.Range("A2").Value = 'radius medial length'
$ws2.Range("B2").Value = 'proximal-distal'
$ws2.Range("C2").Value = '''medialmost part of'' some ''radius'''
$ws2.Range("A3").Value = 'radius diaphysis breadth'
$ws2.Range("B3").Value = 'medial-lateral'
$ws2.Range("C3").Value = 'diaphysis of radius'
$ws2.Range("A4").Value = 'radius diaphysis depth'
$ws2.Range("B4").Value = 'anterior-posterior'
$ws2.Range("C4").Value = 'diaphysis of radius'
$ws2.Range("A5").Value = 'radius proximal articular breadth'
$ws2.Range("B5").Value = 'medial-lateral'
$ws2.Range("A6").Value = 'radius proximal articular depth'
$ws2.Range("B6").Value = 'anterior-posterior'
$ws2.Range("A7").Value = 'radius proximal breadth'
$ws2.Range("B7").Value = 'medial-lateral'
$ws2.Range("C7").Value = 'proximal epiphysis of radius'
$ws2.Range("A8").Value = 'radius distal articular breadth'
$ws2.Range("B8").Value = 'medial-lateral'
$ws2.Range("A9").Value = 'radius distal articular depth'
$ws2.Range("B9").Value = 'anterior-posterior'
$ws2.Range("A10").Value = 'radius distal breadth'
$ws2.Range("B10").Value = 'medial-lateral'
$ws2.Range("C10").Value = 'distal epiphysis of radius'
$ws2.Range("C5").Value = 'humeral facet on radius'
$ws2.Range("C6").Value = 'humeral facet on radius'
$ws2.Range("A11").Value = 'radial condyle breadth'
$ws2.Range("B11").Value = 'medial-lateral'
$ws2.Range("A12").Value = 'ulnar condyle breadth'
$ws2.Range("B12").Value = 'medial-lateral'
$ws2.Range("A13").Value = 'radius lateral length'
$ws2.Range("B13").Value = 'proximal-distal'
$ws2.Range("C13").Value = '''lateral side of'' some radius'
$ws2.Range("A14").Value = 'radius length from the proximal articular surgace to the distal articular surface'
$ws2.Range("B14").Value = 'proximal-distal'

# widen columns A and B on the axis sheet to fit the new content
$ws2.Columns.Item(1).ColumnWidth = 67.1640625
$ws2.Columns.Item(2).ColumnWidth = 15.6640625

# --- back to "trait": pattern-status markers and notes ---
$ws1.Range("B3").Value = 'Y'
$ws1.Range("B4").Value = 'Y'
$ws1.Range("B5").Value = 'Y'

# --- back to "axis" for the remaining "Structure" entries ---
$ws2.Range("C8").Value = 'radio-carpal joint'
$ws2.Range("C9").Value = 'radio-carpal joint'
$ws1.Range("C3").Value = 'radius medial length'
$ws1.Range("C4").Value = 'radius diaphysis breadth'
$ws1.Range("C5").Value = 'radius diaphysis depth'
$ws1.Range("C6").Value = 'radius proximal articular breadth'
$ws1.Range("C7").Value = 'radius proximal articular depth'
$ws1.Range("B8").Value = 'Y'
$ws1.Range("C8").Value = 'radius proximal breadth'
$ws1.Range("C9").Value = 'radius distal articular breadth'
$ws1.Range("C10").Value = 'radius distal articular depth'
$ws1.Range("B11").Value = 'Y'
$ws1.Range("C11").Value = 'radius distal breadth'
$ws1.Range("C12").Value = 'radial condyle breadth'
$ws1.Range("C13").Value = 'ulnar condyle breadth'
$ws1.Range("B14").Value = 'Y'
$ws1.Range("C14").Value = 'radius lateral length'
$ws1.Range("C15").Value = 'radius length from the proximal articular surgace to the distal articular surface'
$ws1.Range("B16").Value = 'Y'
$ws1.Range("C16").Value = 'radius diaphysis circumference'
$ws1.Range("D12").Value = 'lateral distal condyle breadth; radial notch'
$ws1.Range("D13").Value = 'medial distal condyle breadth; ulnar notch'

# --- New sheet "circumference" inserted between "axis" and "structures" ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "circumference"
$ws3.Range("A2").Value = 'radius diaphysis circumference'
$ws3.Range("A2").Select()

# --- selections / active sheet ---
$ws2.Range("C11").Select()
$ws1.Activate()
$ws1.Range("D14").Select()
